$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts all existing data rows down by one,
# row 73 -> row 74, etc.) mirroring the weekly data refresh described by the
# commit message ("Fruta / hortaliza, semanal").
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header row);
# clear that so the new data row matches the plain formatting used by the
# rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Re-apply the date number format used by the other rows in column D.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat()

# Populate the new row 2 with this week's data.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44922
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Provincia de Diguillín"
$ws.Range("P2").Value = 600
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
